# Flow24monthsStatisticsCalculator.xlsx - swap the gauge dataset plugged into the
# "Statistics calculator" sheet from the Smith River / McKenzie gauge to the
# North Santiam / NSantiam gauge (new 24-month flow series), matching the
# author's data refresh. Everything downstream (J:U formulas, averages,
# skill-score cells, the chart) recalculates automatically off these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistics calculator")

# --- Header labels (H3 / I3) describing the new gauge / observation file ---
$ws.Range("H3").Value = " USGS_14178000_flow_NO SANTIAM R BLW BOULDER CRK  NR DETROIT  OR_23780591"
$ws.Range("I3").Value = " Obs:..\Observations\NSantiam\USGS_14178000_flow_NO SANTIAM R BLW BOULDER CRK  NR DETROIT  OR_23780591.csv"

# --- Raw modeled (H) / observed (I) flow values for the 24 months (rows 4-27) ---
$ws.Range("H4").Value = 1318.969116
$ws.Range("I4").Value = 1057.9149170000001
$ws.Range("H5").Value = 654.81872599999997
$ws.Range("I5").Value = 726.88659700000005
$ws.Range("H6").Value = 1309.1667480000001
$ws.Range("I6").Value = 751.378784
$ws.Range("H7").Value = 2291.9421390000002
$ws.Range("I7").Value = 2731.420654
$ws.Range("H8").Value = 1116.070923
$ws.Range("I8").Value = 1012.70343
$ws.Range("H9").Value = 641.40801999999996
$ws.Range("I9").Value = 639.75982699999997
$ws.Range("H10").Value = 464.63085899999999
$ws.Range("I10").Value = 461.15029900000002
$ws.Range("H11").Value = 389.180969
$ws.Range("I11").Value = 373.73623700000002
$ws.Range("H12").Value = 435.92251599999997
$ws.Range("I12").Value = 386.91287199999999
$ws.Range("H13").Value = 543.51409899999999
$ws.Range("I13").Value = 432.39080799999999
$ws.Range("H14").Value = 373.34747299999998
$ws.Range("I14").Value = 381.75524899999999
$ws.Range("H15").Value = 508.02227800000003
$ws.Range("I15").Value = 536.34008800000004
$ws.Range("H16").Value = 1263.4288329999999
$ws.Range("I16").Value = 1170.87085
$ws.Range("H17").Value = 1465.9085689999999
$ws.Range("I17").Value = 1254.3043210000001
$ws.Range("H18").Value = 912.38482699999997
$ws.Range("I18").Value = 690.69152799999995
$ws.Range("H19").Value = 1197.2650149999999
$ws.Range("I19").Value = 1226.988159
$ws.Range("H20").Value = 999.23718299999996
$ws.Range("I20").Value = 1157.6922609999999
$ws.Range("H21").Value = 679.75945999999999
$ws.Range("I21").Value = 827.918274
$ws.Range("H22").Value = 423.25503500000002
$ws.Range("I22").Value = 526.04125999999997
$ws.Range("H23").Value = 356.11215199999998
$ws.Range("I23").Value = 416.39163200000002
$ws.Range("H24").Value = 336.58340500000003
$ws.Range("I24").Value = 381.95684799999998
$ws.Range("H25").Value = 439.83184799999998
$ws.Range("I25").Value = 418.79754600000001
$ws.Range("H26").Value = 772.81073000000004
$ws.Range("I26").Value = 742.23394800000005
$ws.Range("H27").Value = 1337.5745850000001
$ws.Range("I27").Value = 1301.7392580000001

# B4 (NSE score) now needs one extra decimal of display precision, like its
# neighbouring skill-score cells (B6, B12, ...) already have.
$ws.Range("B4").NumberFormat = "0.0000"

# Selection collapses from the old H3:I27 block down to just the header row.
$ws.Range("H3:I3").Select()
